$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marksheet")

# Row 32 - date 45795 - Physics
$ws.Range("C32").Value = "Physics"
$ws.Range("D32").Value = 30
$ws.Range("E32").Value = 28
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 0

# Row 33 - date 45795 - Chemistry
$ws.Range("C33").Value = "Chemistry"
$ws.Range("D33").Value = 35
$ws.Range("E33").Value = 29
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 3

# Row 34 - date 45795 - Biology
$ws.Range("C34").Value = "Biology"
$ws.Range("D34").Value = 50
$ws.Range("E34").Value = 49
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 0

# Row 35 - date 45796 - Physics
$ws.Range("C35").Value = "Physics"
$ws.Range("D35").Value = 25
$ws.Range("E35").Value = 22
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 3

# Row 36 - date 45796 - Chemistry
$ws.Range("C36").Value = "Chemistry"
$ws.Range("D36").Value = 25
$ws.Range("E36").Value = 25
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0

# Row 37 - date 45796 - Biology
$ws.Range("C37").Value = "Biology"
$ws.Range("D37").Value = 50
$ws.Range("E37").Value = 46
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 3

# Row 38 - date 45797 - Physics (subject only, no marks data yet)
$ws.Range("C38").Value = "Physics"

# Row 39 - date 45797 - Chemistry (subject only, no marks data yet)
$ws.Range("C39").Value = "Chemistry"

# Row 40 - date 45797 - Biology
$ws.Range("C40").Value = "Biology"
$ws.Range("D40").Value = 50
$ws.Range("E40").Value = 44
$ws.Range("F40").Value = 3
$ws.Range("G40").Value = 3

# Update sheet view: selection / scroll position
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("D32:D40").Select()

$wb.Save()
